# Auto-generated: apply cryptos.xlsx price/volume refresh (GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.690.30'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '1.596.27'
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'211.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = "'0.0619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = "'19.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.34%  '
$ws.Range("D11").Value = "'0.0837"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.817.06'
$ws.Range("E12").Value = '  -1.64%  '
$ws.Range("D13").Value = '1.594.80'
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("E15").Value = '  -2.75%  '
$ws.Range("D16").Value = "'64.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '26.650.69'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("D19").Value = "'209.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = "'6.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.57%  '
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").Value = "'2.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.98%  '
$ws.Range("D24").Value = "'8.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = "'146.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = "'7.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.60%  '
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").Value = "'15.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = "'0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("D33").Value = "'0.686"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.63%  '
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("D35").Value = '1.295.36'
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("D36").Value = "'2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").Value = "'1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.10%  '
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("D39").Value = "'0.840"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").Value = "'0.791"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = "'2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'5.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = "'63.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '1.730.43'
$ws.Range("E45").Value = '  -1.64%  '
$ws.Range("D46").Value = "'0.902"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.69%  '
$ws.Range("D47").Value = "'89.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0105'
$ws.Range("E49").Value = '  -1.47%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = "'0.0985"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0503"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.48%  '
